$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Value = "Community Services"
[void]$ws.Range("C2").Select()
